$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 16.5
$ws.Range("J2").Value = 5.8
$ws.Range("K2").Value = 6.6
$ws.Range("Y2").Value = 1000
$ws.Range("AD2").Value = 1000
